# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
#
# 1) Reword several subject-line / prompt texts on the "ZansiJourney" sheet.
# 2) Split the "cumcontrol" sheet into two variant sheets:
#       cumcontrol  -> cumcontrol1 (reworded copy of the original content)
#       (new)       -> cumcontrol2 (built from a copy of "dickpic", relabeled
#                       with delay/sync/edge naming + new copy)
#    "dickpic" and "boosters" keep their original content and simply shift
#    position to make room for the newly inserted "cumcontrol2" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. ZansiJourney text tweaks
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("ZansiJourney")
$zj.Range("B8").Value  = "omg"
$zj.Range("B11").Value = "I can't control myself anymore"
$zj.Range("B20").Value = "tell me what you think 😏"
$zj.Range("B22").Value = "wait one sec"

# ---------------------------------------------------------------------------
# 2. Rename "cumcontrol" -> "cumcontrol1" and reword its PPV-control copy
# ---------------------------------------------------------------------------
$cc1 = $wb.Worksheets.Item("cumcontrol")
$cc1.Name = "cumcontrol1"

$cc1.Range("B2").Value = "trust me you want to edge just a little longer for this one"

$cc1.Range("B3").Value = "you're not done until I say you are... open this 😏"
$cc1.Range("C3").Value = "DELAY. Send PPV."

$cc1.Range("B4").Value = "I'm right there too, let's finish this... but you need to see this first"
$cc1.Range("C4").Value = "SYNC variant. Send PPV."

$cc1.Range("B5").Value = "now... right now, with me babe. open this"
$cc1.Range("C5").Value = "SYNC. Send PPV."

$cc1.Range("B6").Value = "not a chance... you're going to wait until I say so"

$cc1.Range("B7").Value = "I didn't say you could cum yet babe 😏"
$cc1.Range("C7").Value = "CONTROL."

# ---------------------------------------------------------------------------
# 3. Duplicate "dickpic" right after "cumcontrol1", rename it "cumcontrol2",
#    and rewrite its rows into the second delay/sync/edge variant set.
#    ("dickpic" itself is left completely untouched and simply slides down
#    to make room for the new sheet; "boosters" likewise shifts position.)
# ---------------------------------------------------------------------------
$dp = $wb.Worksheets.Item("dickpic")
$dp.Copy($null, $cc1)

$cc2 = $wb.Worksheets.Item("dickpic (2)")
$cc2.Name = "cumcontrol2"

$cc2.Range("A2").Value = "delay2"
$cc2.Range("B2").Value = "edge for me... just a little more... this last one is everything"
$cc2.Range("C2").Value = "DELAY variant."

$cc2.Range("A3").Value = "delay1"
$cc2.Range("B3").Value = "hold it... what I'm about to send is the best one and you'll want to last for it 😏"
$cc2.Range("C3").Value = "DELAY. Send PPV."

$cc2.Range("A4").Value = "sync2"
$cc2.Range("B4").Value = "I want us to finish at the same time... this one will push you over"
$cc2.Range("C4").Value = "SYNC variant."

$cc2.Range("A5").Value = "sync1"
$cc2.Range("B5").Value = "okay you earned it babe... let's go together, open this"
$cc2.Range("C5").Value = "SYNC. Send PPV."

$cc2.Range("A6").Value = "edge2"
$cc2.Range("B6").Value = "if you finish without my permission I'll be annoyed"
$cc2.Range("C6").Value = "EDGE variant."

$cc2.Range("A7").Value = "edge1"
$cc2.Range("B7").Value = "slow down babe, I'm in control here 😏"
$cc2.Range("C7").Value = "CONTROL."
